$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newQuestion = "1164. Product Price at a Given Date"
$newDifficulty = "Medium"
$newPattern = "Advanced Select and Joins"
$newLink = "https://leetcode.com/problems/product-price-at-a-given-date/solutions/3825771/simply-using-first-value/?envType=study-plan-v2&envId=top-sql-50 "
$newNotes = "Use first_value(new_price) over window function (partition by product_id order by change_date desc) where change_date <= desired date, union select distinct product_id, where product_id not in (select product_id from Products where change_date <= desired date). You can also use Joins."

# Add the new row (row 30) of data for the new LeetCode problem, matching
# shared-string creation order: Question, then Link, then Notes.
$ws.Range("A30").Value = $newQuestion
$ws.Range("B30").Value = $newDifficulty
$ws.Range("C30").Value = $newPattern
$ws.Range("E30").Value = $newLink
$ws.Range("D30").Value = $newNotes

# Add the hyperlink on the Link cell, then restore the formatting to match
# the rest of the table (copy format down from the row above).
$ws.Hyperlinks.Add($ws.Range("E30"), $newLink) | Out-Null

$ws.Range("B29").Copy()
$ws.Range("B30").PasteSpecial(-4122)
$ws.Range("E29").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Resize the table to include the new row
$table = $ws.ListObjects.Item("Table2")
$table.Resize($ws.Range("A1:E30"))

# Update the active selection to mirror the recorded state
$ws.Range("D36").Select()
